# This script replaces the single autobiographical paragraph in the
# document with the new multi-paragraph "Food Delivery App" write-up,
# preserving the trailing _GoBack bookmark at the very end of the text.

$d = $word.ActiveDocument

# The full original paragraph text (spans several runs / proofErr tags
# in the source markup, but Find/Execute matches across run boundaries).
$oldText = "I am Navya Rajana. I am from Visakhapatnam. I have completed my Tenth Standard in 2012. I have completed My Intermediate College in 2014. I have completed my graduation in B.tech., in felid of Electrical and Electronics Engineering. I am graduated in the year of 2018. I have worked for ChanduSoft Technologies private limited, Visakhapatnam for one year. I was designated as Junior Engineer-ITES at ChanduSoft Technologies private limited. Later I have worked for Amazon Development Centre India Private limited for two years. I was designated as TRON associate at Amazon Development Centre India Private limited. I have taken a time off between ChanduSoft Technologies private limited and Amazon Development Centre India Private limited for nearly one year where I have learned certification course. And after resigning with Amazon Development Centre India Private limited I have joined in a Digital University to learn Data Science. Currently I am learning Data Analysis in the Data Science Course. So far I have achieved certificate of Achievement in Python course. Coming to my hobbies I have keen interest in arts and craft works. In my free time I dedicate myself in bringing out the creativity in me by painting or making new items out of scratch. This is all about myself."

# New content as one big replacement, with "^p" marking paragraph breaks.
# The two genuinely-blank paragraphs use unique placeholder tokens so the
# big replace doesn't collapse them into an (engine-artifact) empty run;
# the placeholders are stripped in a second pass below, leaving true
# empty paragraphs.
$newText = "In Edyoda we have done some projects like Food Delivery App replica, Book my show replica and some other few. Now I am going to describe one of the projects we have done on the platform. I will give detailed script of how we create a Food Delivery App replica.^pFirstly The application will have a log-in for admin and users. Admin will have Some of the functionalities like 1.Add new food items. Food Item will have the following details:^p2.FoodID It should be generated automatically by the application.^p3.Name^p4.Quantity. For eg, 100ml, 250gm, 4pieces etc^p5.Price^p6.Discount^p7.Stock. Amount left in stock in the restaurant.^p8.Edit food items using FoodID.^p9.View the list of all food items.^p10.Remove a food item from the menu using FoodID.^pThe user will have the following functionalities:^p@@BLANK1@@^pRegister on the application. Following to be entered for registration:^p1Full Name^p2Phone Number^p3Email^p4Address^p5Password^p6Log in to the application^p7The user will see 3 options:^p8Place New Order^p9Order History^p10Update Profile^p11Place New Order: The user can place a new order at the restaurant.^p12Show list of food. The list item should as follows:^p1. Tandoori Chicken (4 pieces) [INR 240]^p2. Vegan Burger (1 Piece) [INR 320]^p3. Truffle Cake (500gm) [INR 900]^p13Users should be able to select food by entering an array of numbers. For example, if the user wants to order Vegan Burger and 14Truffle Cake they should enter [2, 3]^p15Once the items are selected user should see the list of all the items selected. The user will also get an option to place an order.^p16Order History should show a list of all the previous orders^p17Update Profile: the user should be able to update their profile.^p@@BLANK2@@^pwith the above mentioned details we need to structure our Food Delivery App replica adding those files in JSON."

# Pass 1: swap in all the new paragraphs in a single Find/Replace.
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Pass 2: clear the blank-paragraph placeholders, leaving clean empty
# paragraphs (no leftover run).
$d.Content.Find.Execute("@@BLANK1@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
$d.Content.Find.Execute("@@BLANK2@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Pass 3: re-add the _GoBack bookmark right after the very last character
# of the document (collapsed / zero-length), matching its original
# position relative to the (now different) final run.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endPos = $lastPara.Range.End - 1

# Workaround: a collapsed range sitting exactly at "end of paragraph text"
# gets mis-resolved to document position 0 by Bookmarks.Add in this
# runtime. Temporarily extend the text by one placeholder character so
# the target position is no longer the paragraph's literal end, add the
# bookmark there, then remove the placeholder again.
$d.Range($endPos, $endPos).InsertAfter("@")
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))
$d.Range($endPos, $endPos + 1).Text = ""
